# [Kadastro App] Yeni kayit eklendi: 2866
#
# Appends a new record row (Kayit No 2866) to both the master "Kayitlar"
# sheet and the district-specific "Erdemli" sheet. The numeric-looking
# fields (Kayit No, Tarih, Parsel Sayisi) are written as literal text so
# they match the existing rows (which store every column as text), not
# as numbers / dates.

$wb = $excel.ActiveWorkbook

function Add-KadastroRow($SheetName, $RowNum, $KayitNo, $Tarih, $Birim, $ParselSayisi, $Is, $Personeller) {
    $ws = $wb.Worksheets.Item($SheetName)

    $rngStr = "A$RowNum" + ":F$RowNum"
    $rng = $ws.Range($rngStr)

    # Force text interpretation so values like "2866" / "2025-09-04" are
    # not auto-converted into numbers / dates.
    $rng.NumberFormat = "@"

    $ws.Range("A$RowNum").Value = $KayitNo
    $ws.Range("B$RowNum").Value = $Tarih
    $ws.Range("C$RowNum").Value = $Birim
    $ws.Range("D$RowNum").Value = $ParselSayisi
    $ws.Range("E$RowNum").Value = $Is
    $ws.Range("F$RowNum").Value = $Personeller

    # Drop the temporary "@" number format again so the new cells end up
    # without any explicit style, same as the surrounding rows.
    $rng.ClearFormats()
}

$kayitNo      = "2866"
$tarih        = "2025-09-04"
$birim        = "Erdemli"
$parselSayisi = "1"
$is           = "MAKS"
$personeller  = "EMİNE ALANLI KIRCILI (K.Mühendisi), AYHAN KARADAYI (K.Teknisyeni)"

# "Kayitlar" currently has data in rows 1-4 -> new row goes to row 5.
Add-KadastroRow "Kayitlar" 5 $kayitNo $tarih $birim $parselSayisi $is $personeller

# "Erdemli" currently has data in rows 1-3 -> new row goes to row 4.
Add-KadastroRow "Erdemli" 4 $kayitNo $tarih $birim $parselSayisi $is $personeller

Write-Host "Added record 2866 to Kayitlar!A5:F5 and Erdemli!A4:F4"
